# Szenario 0 added / weighting rework for "Gewichtung" sheet.
# Updates group_weight (C) and within_group_weight (F) values for every
# criterion row, drops the now-unused total_weight column (G), and moves
# the active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gewichtung")

# --- group_weight (col C) / within_group_weight (col F) updates ---------

# umweltbelastung
$ws.Range("C2:C3").Value = 0.05
$ws.Range("F2").Value = 0.5
$ws.Range("F3").Value = 0.5

# langlebigkeit_wirtschaftlichkeit
$ws.Range("C4:C6").Value = 0.8
$ws.Range("F4").Value = 0.5
$ws.Range("F5").Value = 0.25
$ws.Range("F6").Value = 0.25

# multifunktionale_nutzungsqualitaet
$ws.Range("C7:C10").Value = 0.1
$ws.Range("F7").Value = 0.25
$ws.Range("F8").Value = 0.25
$ws.Range("F9").Value = 0.25
$ws.Range("F10").Value = 0.25

# kreislauffaehigkeit
$ws.Range("C11:C13").Value = 0.05
$ws.Range("F11").Value = 0.333333333
$ws.Range("F12").Value = 0.333333333
$ws.Range("F13").Value = 0.333333333

# --- drop the total_weight column (G) ------------------------------------
$ws.Range("G1:G13").ClearContents()

# --- move the saved selection --------------------------------------------
$ws.Range("D25").Select()
